# Apply "keep highest val strategy": update predicted rating (column D)
# values for the specified rows to match the new predictions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "3" = "A"
    "20" = "BBB"
    "47" = "B"
    "48" = "A"
    "55" = "A"
    "59" = "BB"
    "62" = "A"
    "65" = "BB"
    "69" = "B"
    "71" = "BBB"
    "80" = "BB"
    "81" = "AA"
    "82" = "AA"
    "83" = "AA"
    "84" = "AA"
    "87" = "BBB"
    "94" = "BBB"
    "116" = "BBB"
    "119" = "BBB"
    "120" = "BBB"
    "130" = "BBB"
    "142" = "BBB"
    "144" = "BBB"
    "146" = "BBB"
    "154" = "BBB"
    "159" = "A"
    "166" = "A"
    "167" = "A"
    "182" = "BBB"
    "184" = "B"
    "194" = "BBB"
    "195" = "BBB"
    "203" = "AAA"
    "204" = "AA"
    "207" = "BBB"
    "213" = "BBB"
    "226" = "BBB"
    "229" = "BB"
    "232" = "AA"
    "233" = "BBB"
    "239" = "A"
    "249" = "BB"
    "255" = "BBB"
    "273" = "BBB"
    "290" = "A"
    "293" = "BBB"
    "301" = "A"
    "311" = "BBB"
    "313" = "BB"
    "322" = "BB"
    "323" = "BB"
    "344" = "BB"
    "349" = "BB"
    "367" = "BBB"
    "393" = "BBB"
    "401" = "BBB"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
